$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new family-tree rows 105-150 ("Aisha and Wasfiyya part of the family"),
# following the same layout as the existing rows: A=ROW() formula, B=Sex, C=First Name,
# G/H = Father/Mother lookup formulas (only present where the source data had them).
$rows = @(
  @{ Row=105; Sex='Female'; Name='Aisha'; GFormula='=C$2'; HFormula='=$C$3' }
  @{ Row=106; Sex='Male'; Name='Mohammad'; GFormula=$null; HFormula=$null }
  @{ Row=107; Sex='Male'; Name='Akram'; GFormula='=$C$106'; HFormula='=$C$105' }
  @{ Row=108; Sex='Female'; Name='??'; GFormula=$null; HFormula=$null }
  @{ Row=109; Sex='Male'; Name='Muhammad'; GFormula='=$C$107'; HFormula='=$C$108' }
  @{ Row=110; Sex='Female'; Name='Dima?'; GFormula='=$C$107'; HFormula='=$C$108' }
  @{ Row=111; Sex='Male'; Name='??'; GFormula='=$C$107'; HFormula='=$C$108' }
  @{ Row=112; Sex='Male'; Name='Riyad'; GFormula='=$C$106'; HFormula='=$C$105' }
  @{ Row=113; Sex='Female'; Name='??'; GFormula=$null; HFormula=$null }
  @{ Row=114; Sex='Female'; Name='Mirna'; GFormula='=$C$112'; HFormula='=$C$113' }
  @{ Row=115; Sex='Male'; Name='??'; GFormula=$null; HFormula=$null }
  @{ Row=116; Sex='Male'; Name='??'; GFormula='=$C$115'; HFormula='=$C$114' }
  @{ Row=117; Sex='Male'; Name='??'; GFormula='=$C$115'; HFormula='=$C$114' }
  @{ Row=118; Sex='Female'; Name='Rajai'; GFormula='=$C$112'; HFormula='=$C$113' }
  @{ Row=119; Sex='Male'; Name='??'; GFormula=$null; HFormula=$null }
  @{ Row=120; Sex='Male'; Name='??'; GFormula='=$C$119'; HFormula='=$C$118' }
  @{ Row=121; Sex='Female'; Name='Dalia'; GFormula='=$C$112'; HFormula='=$C$113' }
  @{ Row=122; Sex='Male'; Name='Ziyad'; GFormula='=$C$106'; HFormula='=$C$105' }
  @{ Row=123; Sex='Female'; Name='??'; GFormula=$null; HFormula=$null }
  @{ Row=124; Sex='Male'; Name='Amin'; GFormula='=$C$122'; HFormula='=$C$123' }
  @{ Row=125; Sex='Male'; Name='Murad'; GFormula='=$C$122'; HFormula='=$C$123' }
  @{ Row=126; Sex='Female'; Name='Kamal'; GFormula='=$C$122'; HFormula='=$C$123' }
  @{ Row=127; Sex='Male'; Name='Ali'; GFormula='=$C$106'; HFormula='=$C$105' }
  @{ Row=128; Sex='Female'; Name='Hanija'; GFormula=$null; HFormula=$null }
  @{ Row=129; Sex='Male'; Name='Nasser'; GFormula='=$C$127'; HFormula='=$C$128' }
  @{ Row=130; Sex='Male'; Name='Anwar'; GFormula='=$C$127'; HFormula='=$C$128' }
  @{ Row=131; Sex='Female'; Name='Wasfiyya'; GFormula='=$C$106'; HFormula='=$C$105' }
  @{ Row=132; Sex='Male'; Name='Sadiq'; GFormula=$null; HFormula=$null }
  @{ Row=133; Sex='Female'; Name='Maha'; GFormula='=$C$132'; HFormula='=$C$131' }
  @{ Row=134; Sex='Male'; Name='??'; GFormula=$null; HFormula=$null }
  @{ Row=135; Sex='Male'; Name='Ali'; GFormula='=$C$134'; HFormula='=$C$133' }
  @{ Row=136; Sex='Female'; Name='Raya'; GFormula='=$C$134'; HFormula='=$C$133' }
  @{ Row=137; Sex='Female'; Name='Zaina'; GFormula='=$C$134'; HFormula='=$C$133' }
  @{ Row=138; Sex='Female'; Name='Dana'; GFormula='=$C$134'; HFormula='=$C$133' }
  @{ Row=139; Sex='Female'; Name='Amira'; GFormula='=$C$132'; HFormula='=$C$131' }
  @{ Row=140; Sex='Male'; Name='Maher'; GFormula='=$C$132'; HFormula='=$C$131' }
  @{ Row=141; Sex='Female'; Name='??'; GFormula=$null; HFormula=$null }
  @{ Row=142; Sex='Male'; Name='Sadeq'; GFormula='=$C$140'; HFormula='=$C$141' }
  @{ Row=143; Sex='Male'; Name='??'; GFormula='=$C$140'; HFormula='=$C$141' }
  @{ Row=144; Sex='Female'; Name='Randa'; GFormula='=$C$132'; HFormula='=$C$131' }
  @{ Row=145; Sex='Female'; Name='Dalal'; GFormula='=$C$132'; HFormula='=$C$131' }
  @{ Row=146; Sex='Male'; Name='??'; GFormula=$null; HFormula=$null }
  @{ Row=147; Sex='Male'; Name='??'; GFormula='=$C$146'; HFormula='=$C$145' }
  @{ Row=148; Sex='Male'; Name='??'; GFormula='=$C$146'; HFormula='=$C$145' }
  @{ Row=149; Sex='Male'; Name='??'; GFormula='=$C$146'; HFormula='=$C$145' }
  @{ Row=150; Sex='Female'; Name='Rima'; GFormula='=$C$132'; HFormula='=$C$131' }
)
foreach ($r in $rows) {
  $rowNum = $r.Row

  # Column A: running-row-number formula, same as every other data row.
  $ws.Range("A$rowNum").Formula = "=ROW()"

  # Column B: Sex value + the matching fill style (copy format from an existing
  # Male/Female reference cell so the style index is reused, not duplicated).
  if ($r.Sex -eq 'Male') {
    $ws.Range("B2").Copy()
  } else {
    $ws.Range("B3").Copy()
  }
  $ws.Range("B$rowNum").PasteSpecial(-4122)
  $ws.Range("B$rowNum").Value = $r.Sex

  # Column C: First name.
  $ws.Range("C$rowNum").Value = $r.Name

  # Columns G/H: Father / Mother formulas, when present in the source row.
  if ($r.GFormula) {
    $ws.Range("G$rowNum").Formula = $r.GFormula
  }
  if ($r.HFormula) {
    $ws.Range("H$rowNum").Formula = $r.HFormula
  }
}

$excel.CutCopyMode = 0

# Move the selection the same way the author ended up after entering the
# new rows (bottom of the newly-added data, column C).
$ws.Range("A124").Select()
$ws.Range("C150").Select()

Write-Output "Added rows 105-150 ($($rows.Count) people)."
